# TypeScript SDK facade API - basics (#78)
# Adds two new worksheets ("Helloworld!" and "Doesitwork "), adds a
# "Bool:" row with a boolean value to Sheet1, and tweaks a couple of
# selections/styles along the way.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)

# "Helloworld!" is created first (so it is assigned the lower internal
# sheetId) and lives at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$helloSheet = $wb.Worksheets.Add($null, $lastSheet)
$helloSheet.Name = "Helloworld!"

# "Doesitwork " is created second, but placed right after Sheet1 (i.e.
# before "Helloworld!" in tab order).
$doesSheet = $wb.Worksheets.Add($null, $sheet1)
$doesSheet.Name = "Doesitwork "

# Populate the new "Doesitwork " sheet.
$doesSheet.Range("A1").Value = "D"
[void]$doesSheet.Range("C6").Select()

# Back on Sheet1: add a boolean row under the existing data.
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A4").Value = "Bool:"
$ws.Range("B4").Value = $true

# B2 drops its custom number-format style, going back to the default.
$ws.Range("B2").Style = "Normal"

# Move the active selection on Sheet1.
[void]$ws.Range("B7").Select()
